$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts old rows 8-14 down to 9-15)
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the new reference entry.
# Order matches the shared-string table insertion order observed in the
# target workbook (title, link, author) rather than column order.
$ws.Cells.Item(8, 2).Value = "Forecasting: Principles and Practice, the Pythonic Way"
$ws.Cells.Item(8, 3).Value = "https://otexts.com/fpppy/"
$ws.Cells.Item(8, 1).Value = "Hyndman, R.J., Athanasopoulos, G., Garza, A., Challu, C., Mergenthaler, M., & Olivares, K.G."

# Match formatting used by the rest of the table (wrap text style) and row height
$ws.Rows.Item(8).RowHeight = 51
$ws.Range("A8:C8").WrapText = $true

# Update view state to match the target: active selection moves to C8
$ws.Activate()
$ws.Range("C8").Select()
